# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet
#    that carries it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) Narrow the "Status" column(s) that used to hold that long string:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status values -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Shrink the now-narrower status columns ------------------------------------
# Target stored column width is ~13.41 (down from ~17.22); the engine stores
# ColumnWidth quantized in character units, so feed it the character width
# that lands on the closest quantized bucket to the target.
$newStatusWidth = 12.5

$overview.Range("E1:F1").ColumnWidth = $newStatusWidth
$zhcn.Range("C1").ColumnWidth = $newStatusWidth
$dede.Range("C1").ColumnWidth = $newStatusWidth
